# correção nos dados e inicio da analise PNAD 2009
#
# The row for "grandes regiões e unidades da federação" (row 6) was an
# orphan header-like label with no data beneath it on its own row; it is
# removed entirely, and every row below it (the actual "norte" data and
# all the following regions/states) shifts up by one row. This also drops
# the now-unused shared string and removes the trailing blank slot that
# used to be occupied by the last row (distrito federal), which now lands
# on row 37 instead of row 38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the entire row shifts rows 7:38 up to 6:37, keeping each
# region's data intact (norte's values move from row 7 to row 6, etc.)
# and Excel drops the shared string that's no longer referenced by any
# cell when the file is saved.
$ws.Rows("6:6").Delete()
